$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("E2").Value = "Motul Oil 2.jpg,Motul Oil.jpg"
$ws.Range("E2").Select()
